$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "DEND" = "Dendropsophus_microcephalus"
    "TRGL" = "Troglodytes_aedon"
    "FRAG" = "Leptodactylus_fragilis"
    "FUSC" = "Leptodactylus_fuscus"
    "BOAN" = "Boana_platanera"
    "MMLS" = "Alouatta_sp"
    "PTGN" = "Patagioenas_cayennensis"
    "DDPX" = "Dendroplex_picus"
}

for ($r = 2; $r -le 109; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $code = $cell.Value2
    if ($map.ContainsKey($code)) {
        $cell.Value = $map[$code]
    }
}
